$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# Collapse the opening of the Token Ring paragraph (which previously said
# it is an OSI-layer-2 Data-Link protocol / competitor to Ethernet, then
# repeated "Like Ethernet, Token Ring utilises...") down to a single,
# tighter opening sentence. The trailing bold sentence ("Token Ring is an
# implementation of point-to-point communication.") is left untouched.
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Replacement.ClearFormatting()
$find1.Text = "Within the Open Systems Interconnection (OSI) model, Token Ring is a layer 2 Data-Link layer protocol and is a direct competitor to Ethernet. Like Ethernet, Token Ring utilises "
$find1.Replacement.Text = "Token Ring networks utilise "
$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2) | Out-Null

# --- Edit 2 -----------------------------------------------------------
# Add a new sentence right after "...all stations are connected
# sequentially." explaining that each station connects to two neighbours.
$rng2 = $d.Content
$find2 = $rng2.Find
$find2.ClearFormatting()
$find2.Text = "all stations are connected sequentially."
$find2.Execute() | Out-Null
$rng2.Collapse(0)
$rng2.InsertAfter(" Therefore, a each station is connected to a total of two other stations, the previous and next station in the network.")
